$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.614.43"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "3.250.03"
$ws.Range("E3").Value = "  +2.79%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "607.49"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.249.97"
$ws.Range("E8").Value = "  +2.86%  "
$ws.Range("D9").Value = "0.549"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("E10").Value = "  +2.31%  "
$ws.Range("E11").Value = "  +5.36%  "
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("E13").Value = "  +2.34%  "
$ws.Range("D14").Value = "39.14"
$ws.Range("E14").Value = "  +2.18%  "
$ws.Range("D15").Value = "3.789.29"
$ws.Range("E15").Value = "  +2.98%  "
$ws.Range("D16").Value = "66.659.58"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").Value = "3.255.23"
$ws.Range("E18").Value = "  +3.05%  "
$ws.Range("D19").Value = "0.113"
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("D20").Value = "508.68"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "15.41"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").Value = "0.750"
$ws.Range("E22").Value = "  +3.27%  "
$ws.Range("D23").Value = "8.11"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("D24").Value = "14.81"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").Value = "86.77"
$ws.Range("E25").Value = "  +2.74%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").Value = "0.140"
$ws.Range("E27").Value = "  +58.62%  "
$ws.Range("E28").Value = "  +1.75%  "
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("D31").Value = "2.90"
$ws.Range("E31").Value = "  -4.92%  "
$ws.Range("D32").Value = "6.84"
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("D33").Value = "28.10"
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("D35").Value = "1.16"
$ws.Range("E35").Value = "  -3.19%  "
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("E37").Value = "  +23.02%  "
$ws.Range("D38").Value = "0.0₃0797"
$ws.Range("E38").Value = "  +18.52%  "
$ws.Range("D39").Value = "55.64"
$ws.Range("E39").Value = "  +1.62%  "
$ws.Range("D40").Value = "494.54"
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("D41").Value = "0.0426"
$ws.Range("E41").Value = "  +1.98%  "
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("D43").Value = "8.83"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").Value = "0.293"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").Value = "2.52"
$ws.Range("E45").Value = "  +4.18%  "
$ws.Range("D46").Value = "2.971.06"
$ws.Range("E46").Value = "  +5.46%  "
$ws.Range("D47").Value = "28.90"
$ws.Range("E47").Value = "  +3.94%  "
$ws.Range("E48").Value = "  +5.28%  "
$ws.Range("D49").Value = "0.119"
$ws.Range("E49").Value = "  +2.81%  "
$ws.Range("D51").Value = "121.36"
$ws.Range("E51").Value = "  +0.69%  "
